# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# D-column values that would otherwise be auto-parsed by Excel as numbers
# (losing significant trailing/leading zeros, e.g. "19.30" -> 19.3) are
# forced to Text format first so the literal string is preserved exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.424.25"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.105.05"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.22"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5231"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4571"
$ws.Range("E8").Value = "  +6.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.14"
$ws.Range("E9").Value = "  +16.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08918"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.177"
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.42"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "2.091.80"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.796"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.017"
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.42"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.30"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.329"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").Value = "30.476.52"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.365"
$ws.Range("E25").Value = "  +3.33%  "
$ws.Range("D26").Value = "2.341.23"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.34"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.573"
$ws.Range("E28").Value = "  +3.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.72"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.57"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.231"
$ws.Range("E31").Value = "  +5.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.743"
$ws.Range("E32").Value = "  +17.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1072"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.192"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.930"
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.52"
$ws.Range("E36").Value = "  +10.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02578"
$ws.Range("E37").Value = "  +1.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06833"
$ws.Range("E38").Value = "  +4.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.550"
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.81"
$ws.Range("E40").Value = "  +3.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2303"
$ws.Range("E41").Value = "  +3.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6898"
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.349"
$ws.Range("E44").Value = "  +7.91%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6376"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.96"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.662"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000352"
$ws.Range("E49").Value = "  +25.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.247"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.55"
